$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.3232956666666666
$ws.Range("H2").Value = 0.9698869999999999
$ws.Range("I2").Value = 0.09891704828958615
$ws.Range("J2").Value = 0.09891704828958615
$ws.Range("M2").Value = 2.245342666666666
$ws.Range("N2").Value = 6.736027999999999
$ws.Range("O2").Value = 0.02318640424300622
$ws.Range("P2").Value = 0.02318640424300622
$ws.Range("Q2").Value = 0.725909554315111
$ws.Range("R2").Value = 6.533185988835999
$ws.Range("S2").Value = 0.002293530668167311
$ws.Range("T2").Value = 0.002293530668167312
# Row 3
$ws.Range("G3").Value = 0.3232956666666666
$ws.Range("H3").Value = 0.9698869999999999
$ws.Range("I3").Value = 0.09891704828958615
$ws.Range("J3").Value = 0.09891704828958615
$ws.Range("O3").Value = 0.225628233631131
$ws.Range("P3").Value = 0.225628233631131
$ws.Range("Q3").Value = 7.063867635512444
$ws.Range("R3").Value = 63.57480871961199
$ws.Range("S3").Value = 0.02231847888158461
$ws.Range("T3").Value = 0.02231847888158461
# Row 4
$ws.Range("G4").Value = 0.3232956666666666
$ws.Range("H4").Value = 0.9698869999999999
$ws.Range("I4").Value = 0.09891704828958615
$ws.Range("J4").Value = 0.09891704828958615
$ws.Range("M4").Value = 72.46803266666666
$ws.Range("N4").Value = 217.404098
$ws.Range("O4").Value = 0.7483370467453728
$ws.Range("P4").Value = 0.7483370467453727
$ws.Range("Q4").Value = 23.42860093299177
$ws.Range("R4").Value = 210.857408396926
$ws.Range("S4").Value = 0.07402329178979833
$ws.Range("T4").Value = 0.07402329178979833
# Row 5
$ws.Range("G5").Value = 0.3232956666666666
$ws.Range("H5").Value = 0.9698869999999999
$ws.Range("I5").Value = 0.09891704828958615
$ws.Range("J5").Value = 0.09891704828958615
$ws.Range("M5").Value = 0.2758273333333334
$ws.Range("N5").Value = 0.827482
$ws.Range("O5").Value = 0.002848315380489998
$ws.Range("P5").Value = 0.002848315380489997
$ws.Range("Q5").Value = 0.0891737816148889
$ws.Range("R5").Value = 0.802564034534
$ws.Range("S5").Value = 0.0002817469500359001
$ws.Range("T5").Value = 0.0002817469500359
# Row 6
$ws.Range("I6").Value = 0.6064896735907829
$ws.Range("J6").Value = 0.6064896735907829
$ws.Range("M6").Value = 2.245342666666666
$ws.Range("N6").Value = 6.736027999999999
$ws.Range("O6").Value = 0.02318640424300622
$ws.Range("P6").Value = 0.02318640424300622
$ws.Range("Q6").Value = 4.450766134510221
$ws.Range("R6").Value = 40.056895210592
$ws.Range("S6").Value = 0.01406231474108478
$ws.Range("T6").Value = 0.01406231474108478
# Row 7
$ws.Range("I7").Value = 0.6064896735907829
$ws.Range("J7").Value = 0.6064896735907829
$ws.Range("O7").Value = 0.225628233631131
$ws.Range("P7").Value = 0.225628233631131
$ws.Range("S7").Value = 0.1368411937678095
$ws.Range("T7").Value = 0.1368411937678095
# Row 8
$ws.Range("I8").Value = 0.6064896735907829
$ws.Range("J8").Value = 0.6064896735907829
$ws.Range("M8").Value = 72.46803266666666
$ws.Range("N8").Value = 217.404098
$ws.Range("O8").Value = 0.7483370467453728
$ws.Range("P8").Value = 0.7483370467453727
$ws.Range("Q8").Value = 143.6476803365636
$ws.Range("R8").Value = 1292.829123029072
$ws.Range("S8").Value = 0.4538586912164916
$ws.Range("T8").Value = 0.4538586912164915
# Row 9
$ws.Range("I9").Value = 0.6064896735907829
$ws.Range("J9").Value = 0.6064896735907829
$ws.Range("M9").Value = 0.2758273333333334
$ws.Range("N9").Value = 0.827482
$ws.Range("O9").Value = 0.002848315380489998
$ws.Range("P9").Value = 0.002848315380489997
$ws.Range("Q9").Value = 0.5467508244497779
$ws.Range("R9").Value = 4.920757420048001
$ws.Range("S9").Value = 0.001727473865396985
$ws.Range("T9").Value = 0.001727473865396985
# Row 10
$ws.Range("G10").Value = 0.6996333333333333
$ws.Range("H10").Value = 2.0989
$ws.Range("I10").Value = 0.2140630740024481
$ws.Range("J10").Value = 0.2140630740024481
$ws.Range("M10").Value = 2.245342666666666
$ws.Range("N10").Value = 6.736027999999999
$ws.Range("O10").Value = 0.02318640424300622
$ws.Range("P10").Value = 0.02318640424300622
$ws.Range("Q10").Value = 1.570916574355555
$ws.Range("R10").Value = 14.1382491692
$ws.Range("S10").Value = 0.004963352967321317
$ws.Range("T10").Value = 0.004963352967321318
# Row 11
$ws.Range("G11").Value = 0.6996333333333333
$ws.Range("H11").Value = 2.0989
$ws.Range("I11").Value = 0.2140630740024481
$ws.Range("J11").Value = 0.2140630740024481
$ws.Range("O11").Value = 0.225628233631131
$ws.Range("P11").Value = 0.225628233631131
$ws.Range("Q11").Value = 15.28667956182222
$ws.Range("R11").Value = 137.5801160564
$ws.Range("S11").Value = 0.04829867327282245
$ws.Range("T11").Value = 0.04829867327282245
# Row 12
$ws.Range("G12").Value = 0.6996333333333333
$ws.Range("H12").Value = 2.0989
$ws.Range("I12").Value = 0.2140630740024481
$ws.Range("J12").Value = 0.2140630740024481
$ws.Range("M12").Value = 72.46803266666666
$ws.Range("N12").Value = 217.404098
$ws.Range("O12").Value = 0.7483370467453728
$ws.Range("P12").Value = 0.7483370467453727
$ws.Range("Q12").Value = 50.70105125468888
$ws.Range("R12").Value = 456.3094612922
$ws.Range("S12").Value = 0.1601913286162282
$ws.Range("T12").Value = 0.1601913286162282
# Row 13
$ws.Range("G13").Value = 0.6996333333333333
$ws.Range("H13").Value = 2.0989
$ws.Range("I13").Value = 0.2140630740024481
$ws.Range("J13").Value = 0.2140630740024481
$ws.Range("M13").Value = 0.2758273333333334
$ws.Range("N13").Value = 0.827482
$ws.Range("O13").Value = 0.002848315380489998
$ws.Range("P13").Value = 0.002848315380489997
$ws.Range("Q13").Value = 0.1929779966444445
$ws.Range("R13").Value = 1.7368019698
$ws.Range("S13").Value = 0.0006097191460761415
$ws.Range("T13").Value = 0.0006097191460761414
# Row 14
$ws.Range("E14").Value = 3.0
$ws.Range("F14").Value = 1.0
$ws.Range("G14").Value = 0.263201
$ws.Range("H14").Value = 0.789603
$ws.Range("I14").Value = 0.08053020411718284
$ws.Range("J14").Value = 0.08053020411718284
$ws.Range("M14").Value = 2.245342666666666
$ws.Range("N14").Value = 6.736027999999999
$ws.Range("O14").Value = 0.02318640424300622
$ws.Range("P14").Value = 0.02318640424300622
$ws.Range("Q14").Value = 0.5909764352093333
$ws.Range("R14").Value = 5.318787916884
$ws.Range("S14").Value = 0.001867205866432805
$ws.Range("T14").Value = 0.001867205866432805
# Row 15
$ws.Range("E15").Value = 3.0
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = 0.263201
$ws.Range("H15").Value = 0.789603
$ws.Range("I15").Value = 0.08053020411718284
$ws.Range("J15").Value = 0.08053020411718284
$ws.Range("O15").Value = 0.225628233631131
$ws.Range("P15").Value = 0.225628233631131
$ws.Range("Q15").Value = 5.750825690625334
$ws.Range("R15").Value = 51.757431215628
$ws.Range("S15").Value = 0.0181698877089144
$ws.Range("T15").Value = 0.0181698877089144
# Row 16
$ws.Range("E16").Value = 3.0
$ws.Range("F16").Value = 1.0
$ws.Range("G16").Value = 0.263201
$ws.Range("H16").Value = 0.789603
$ws.Range("I16").Value = 0.08053020411718284
$ws.Range("J16").Value = 0.08053020411718284
$ws.Range("M16").Value = 72.46803266666666
$ws.Range("N16").Value = 217.404098
$ws.Range("O16").Value = 0.7483370467453728
$ws.Range("P16").Value = 0.7483370467453727
$ws.Range("Q16").Value = 19.07365866589933
$ws.Range("R16").Value = 171.662927993094
$ws.Range("S16").Value = 0.06026373512285467
$ws.Range("T16").Value = 0.06026373512285466
# Row 17
$ws.Range("E17").Value = 3.0
$ws.Range("F17").Value = 1.0
$ws.Range("G17").Value = 0.263201
$ws.Range("H17").Value = 0.789603
$ws.Range("I17").Value = 0.08053020411718284
$ws.Range("J17").Value = 0.08053020411718284
$ws.Range("M17").Value = 0.2758273333333334
$ws.Range("N17").Value = 0.827482
$ws.Range("O17").Value = 0.002848315380489998
$ws.Range("P17").Value = 0.002848315380489997
$ws.Range("Q17").Value = 0.07259802996066668
$ws.Range("R17").Value = 0.6533822696460001
$ws.Range("S17").Value = 0.0002293754189809708
$ws.Range("T17").Value = 0.0002293754189809708
